$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) text updates
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "機械可読ドキュメント"
$ws.Range("H1").Value = "年"
$ws.Range("I1").Value = "IIIFマニフェストURI"
$ws.Range("J1").Value = "viewingDirection"
$ws.Range("K1").Value = "帰属"
$ws.Range("M1").Value = "ソート用項目"
$ws.Range("N1").Value = "コレクション"
$ws.Range("O1").Value = "サムネイル"
$ws.Range("P1").Value = "ウェブサイトURL"

# ---------------------------------------------------------------------------
# 2. Data row (row 2) text updates
# ---------------------------------------------------------------------------
# H2 becomes empty (its old thumbnail-link content moves to O2)
$ws.Range("H2").ClearContents()
# J2 gains the rightToLeftDirection link text (used to live, un-linked, nowhere)
$ws.Range("J2").Value = "http://iiif.io/api/presentation/2#rightToLeftDirection"
# K2 becomes the attribution text (used to live in M2)
$ws.Range("K2").Value = "東京大学総合図書館 General Library in the University of Tokyo, JAPAN"
# M2 becomes empty (its old attribution content moves to K2)
$ws.Range("M2").ClearContents()
# N2 gains the title text (used to live, un-linked, in K2)
$ws.Range("N2").Value = "水野家古文書(水野忠幹氏旧蔵書文書)"
# O2 becomes the thumbnail image link (used to live in H2)
$ws.Range("O2").Value = "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/30/full/200,151/0/default.jpg"

# ---------------------------------------------------------------------------
# 3. Formatting: move the "hyperlink" cell style (underline + blue font) off
#    H2 (no longer a link) and onto J2 (now a link). The other styled cells
#    (C2, D2, F2, I2, O2, P2) keep whatever formatting they already have.
# ---------------------------------------------------------------------------
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("H2").ClearFormats()

# ---------------------------------------------------------------------------
# 4. Hyperlinks: rebuild from scratch so relationship ids come out in the
#    exact order the target file expects (rId1..rId7 in document order).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/mizuno/document/d28ad412-0d21-4a16-ab81-840fe3ffde71")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/12")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/12/manifest")
$ws.Hyperlinks.Add($ws.Range("J2"), "http://iiif.io/api/presentation/2", "rightToLeftDirection")
$ws.Hyperlinks.Add($ws.Range("O2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/30/full/200,151/0/default.jpg")
$ws.Hyperlinks.Add($ws.Range("P2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/mizuno/")

Write-Output "done"
